$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 12:52"

# --- Countries re-ranked by total cases: Afganistan overtook Cuba/Oman/Uruguay ---
# (rows 88-91 shift down by one position, Afganistan now leads the group)
$ws.Cells.Item(88, 1).Value = "Afganistan"
$ws.Cells.Item(89, 1).Value = "Cuba"
$ws.Cells.Item(90, 1).Value = "Oman"
$ws.Cells.Item(91, 1).Value = "Uruguay"

# --- Countries re-ranked by total cases: Malta overtook Ghana/San Marino ---
$ws.Cells.Item(100, 1).Value = "Malta"
$ws.Cells.Item(101, 1).Value = "Ghana"
$ws.Cells.Item(102, 1).Value = "San Marino"

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Refresh numeric stats (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) for the affected rows ---

Set-Row 14  23514 234 9800  12801 391 18 913   # Suiza
Set-Row 48  2376  166 206   2164  37  0  6     # Catar
Set-Row 83  582   6   62    501   29  0  19    # Libano
Set-Row 88  484   40  32    437   0   1  15    # now Afganistan
Set-Row 89  457   0   27    418   15  0  12    # now Cuba
Set-Row 90  457   38  109   346   3   0  2     # now Oman
Set-Row 91  456   0   192   257   14  0  7     # now Uruguay
Set-Row 100 337   38  16    319   4   1  2     # now Malta
Set-Row 101 313   0   34    273   2   0  6     # now Ghana
Set-Row 102 308   0   45    229   14  0  34    # now San Marino
Set-Row 141 56    1   4     50    2   0  2     # Etiopia
